# Actualización automática de tasas-transfi.xlsx

$wb = $excel.ActiveWorkbook

# --- Hoja1: update the daily conversion note text ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$wsHoja1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 7.54 = 30708.11 pesos`n✅ 30708.11 pesos = 7.54 = 960.38 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- tasas: update the rate cells ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 132.62
$wsTasas.Range("O10").Value = 4072.51
$wsTasas.Range("N12").Value = 4072
$wsTasas.Range("O12").Value = 127.35
